# Diary entry for week 9 (rows 57-65, dates 2020-02-27 .. 2020-03-04)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Dates (column A) ---
$ws.Range("A57").Value = 43888
$ws.Range("A59").Value = 43890
$ws.Range("A61").Value = 43891
$ws.Range("A63").Value = 43892
$ws.Range("A65").Value = 43894

# Give the new date cells the same left-aligned date number format already
# used by the other date cells in column A (copy format from A51).
$ws.Range("A51").Copy()
$ws.Range("A57").PasteSpecial(-4122)
$ws.Range("A59").PasteSpecial(-4122)
$ws.Range("A61").PasteSpecial(-4122)
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("A67").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 57: 2/27 - Learn about design patterns ---
$ws.Range("B57").Value = "5pm-7.50pm"
$ws.Range("C57").Value = "N/A"
$ws.Range("D57").Value = "Learn about design patterns"
$ws.Range("E57").Value = "Learned 3 more key expert practices and design patterns`nIt was really inspiring to hear the experiences of Alberto"
$ws.Range("F57").Value = "Understood how we can make use of design patterns to effectively solve a software design problem and how it helps to structure the codebase in a better way which in turn can improve the readbility of the code"
$ws.Range("G57").Value = "Wish we could spend more time on design patterns and do some practice on that"

# --- Row 59: 2/29 - Discuss this week's homework ---
$ws.Range("D59").Value = "Discuss this week's homework"
$ws.Range("B59").Value = "5pm-7pm"
$ws.Range("C59").Value = "Aman, Vaishakhi"
$ws.Range("G59").Value = "Need to find more simpler issues for the  first issue to fix"

# --- Row 61: 3/1 - Find 5 design patterns in the code ---
$ws.Range("F61").Value = "Understood the importance of having proper naming conventions for classes as it helped in narrowing down the design patterns used. For example there was a class names `"SessionFactory`" which was a clear indicator of the Factory pattern"
$ws.Range("B61").Value = "5pm-9pm"
$ws.Range("C61").Value = "Aman, Vaishakhi"

# --- Row 63: 3/2 - Identify a simple issue and fix it ---
$ws.Range("D63").Value = "Identify a simple issue and fix it"
$ws.Range("B63").Value = "9pm-12pm"
$ws.Range("C63").Value = "Aman, Vaishakhi"
$ws.Range("G63").Value = "Happy to have narrowed down an issue for the first assignment that everyone understood"
$ws.Range("E63").Value = "Looked for couple of more simple issues to fix . Found one issue which all of us could understand properly. Started making code changes"

# --- Row 65: 3/4 - Submit the pull request ---
$ws.Range("D65").Value = "Submit the pull request"
$ws.Range("B65").Value = "9pm-12pm"
$ws.Range("C65").Value = "Aman, Vaishakhi"
$ws.Range("F59").Value = "We realized that certain issues that we thought to be easy were actually not that easy to fix due to either lack of proper test cases or complex dependencies in the codebase. "
$ws.Range("F63").Value = "The issue was easier to understand because 1)  it was a relatively simple change and had lesser dependencies  2) comments in the pull request were really helpful."
$ws.Range("F65").Value = "Coding the issue was not difficult as we had already discussed the scenarios properly and noted down the various conditions that needs to be covered"
$ws.Range("G65").Value = "Hoping that our pull request will get accepted"
$ws.Range("E65").Value = "Completed the code changes, tested it and submitted the pull request`nFinished the report"
$ws.Range("D61").Value = "Find 5 design patterns in the code"
$ws.Range("E61").Value = "Tried to understand various design patterns by going through websites and YouTube videos. Identified 5 design patterns in the codebase and documented the same - Factory, Singleton, Adaptor, Decorator and Command "
$ws.Range("G61").Value = "Happy to have identified 5 patterns"
$ws.Range("E59").Value = "As Aman and Vaishakhi were not present in the last class, I explained this week's homework. We also discussed couple of possible issues to fix."

# --- Row heights to match the wrapped text now filling these rows ---
$ws.Rows(57).RowHeight = 109.2
$ws.Rows(59).RowHeight = 93.6
$ws.Rows(61).RowHeight = 109.2
$ws.Rows(63).RowHeight = 78
$ws.Rows(65).RowHeight = 78

# --- View state: scrolled down to this week's entries, zoomed to 89% ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 89
$ws.Range("E62").Select()
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1
